# Update "想去人数" (F) and "最低票价" (G) figures on the "展览" and
# "全部类型" sheets to the newly scraped values.

$wb = $excel.ActiveWorkbook

# Row => hashtable of column => new value
$updates = @{
    2  = @{ F = 109 }
    4  = @{ F = 1381 }
    5  = @{ F = 1614; G = 50 }
    6  = @{ F = 355 }
    7  = @{ F = 466 }
    9  = @{ F = 193 }
    12 = @{ F = 123 }
    16 = @{ F = 1790 }
    20 = @{ F = 707 }
    23 = @{ F = 4327 }
    25 = @{ F = 304 }
    26 = @{ F = 1146 }
    27 = @{ F = 501 }
    29 = @{ F = 683 }
    31 = @{ F = 342 }
    33 = @{ F = 177 }
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $cols = $updates[$row]
        foreach ($col in $cols.Keys) {
            $ws.Range("$col$row").Value = $cols[$col]
        }
    }
}
